# Add a new final slide ("Lähteet" / Sources) by duplicating the last
# existing slide (Slide 8, "Vertauskuva") so that it inherits the same
# layout, placeholder structure, and formatting, then replace its title
# text.

$p = $ppt.ActivePresentation

$lastIndex = $p.Slides.Count
$sourceSlide = $p.Slides.Item($lastIndex)

$newRange = $sourceSlide.Duplicate()
$newSlide = $newRange.Item(1)

$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Lähteet"
